# Apply update for 2022-11-01 data add:
# - rename sheet/tab to reflect new date
# - update header label "2022 (through 10-23)" -> "2022 (through 10-24)"
# - update I11 85 -> 89
# - update I14 1362 -> 1366

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-10-24"
$ws.Range("I1").Value = "2022 (through 10-24)"
$ws.Range("I11").Value = 89
$ws.Range("I14").Value = 1366
